$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "319÷9=35, 4"
$t.Cell(1, 2).Range.Text = "659÷5=131, 4"
$t.Cell(1, 3).Range.Text = "832÷9=92, 4"
$t.Cell(1, 4).Range.Text = "816÷5=163, 1"
$t.Cell(1, 5).Range.Text = "350÷2=175, 0"

$t.Cell(5, 1).Range.Text = "936÷9=104, 0"
$t.Cell(5, 2).Range.Text = "216÷6=36, 0"
$t.Cell(5, 3).Range.Text = "666÷9=74, 0"
$t.Cell(5, 4).Range.Text = "426÷3=142, 0"
$t.Cell(5, 5).Range.Text = "982÷5=196, 2"

$t.Cell(9, 1).Range.Text = "817÷7=116, 5"
$t.Cell(9, 2).Range.Text = "504÷6=84, 0"
$t.Cell(9, 3).Range.Text = "109÷2=54, 1"
$t.Cell(9, 4).Range.Text = "200÷4=50, 0"
$t.Cell(9, 5).Range.Text = "656÷7=93, 5"

$t.Cell(13, 1).Range.Text = "513÷7=73, 2"
$t.Cell(13, 2).Range.Text = "905÷6=150, 5"
$t.Cell(13, 3).Range.Text = "829÷2=414, 1"
$t.Cell(13, 4).Range.Text = "830÷5=166, 0"
$t.Cell(13, 5).Range.Text = "258÷4=64, 2"

$t.Cell(17, 1).Range.Text = "727÷4=181, 3"
$t.Cell(17, 2).Range.Text = "806÷3=268, 2"
$t.Cell(17, 3).Range.Text = "155÷4=38, 3"
$t.Cell(17, 4).Range.Text = "875÷6=145, 5"
$t.Cell(17, 5).Range.Text = "796÷3=265, 1"

